# Prefix the document's Heading1/2/3 paragraphs with their outline numbers
# (e.g. "Introduction" -> "1. Introduction", "Definitions" -> "1.1 Definitions").
# Matching is done on BOTH the paragraph style and the exact current heading
# text, so body-text runs that happen to repeat a heading's words (e.g. the
# bold "Safe projects" lead-in sentence, or "Process" inside the title) are
# left untouched.

$d = $word.ActiveDocument

$renames = @(
    @{ Style = "Heading 1"; Old = "Introduction"; New = "1. Introduction" },
    @{ Style = "Heading 2"; Old = "Definitions"; New = "1.1 Definitions" },
    @{ Style = "Heading 2"; Old = "Key Principles"; New = "1.2 Key Principles" },
    @{ Style = "Heading 2"; Old = "Transparency reporting"; New = "1.3 Transparency reporting" },
    @{ Style = "Heading 1"; Old = "ODAP Review processes"; New = "2. ODAP Review processes" },
    @{ Style = "Heading 2"; Old = "Simplified review process"; New = "2.1 Simplified review process" },
    @{ Style = "Heading 2"; Old = "Dataset availability"; New = "2.2 Dataset availability" },
    @{ Style = "Heading 2"; Old = "Accreditation of Researchers"; New = "2.3 Accreditation of Researchers" },
    @{ Style = "Heading 3"; Old = "Review of researcher credentials by ODAP data access team"; New = "2.3.1 Review of researcher credentials by ODAP data access team" },
    @{ Style = "Heading 3"; Old = "Decision and notification"; New = "2.3.2 Decision and notification" },
    @{ Style = "Heading 2"; Old = "Bona Fide researchers"; New = "2.4 Bona Fide researchers" },
    @{ Style = "Heading 1"; Old = "Review of Research Proposals"; New = "3. Review of Research Proposals" },
    @{ Style = "Heading 2"; Old = "Safe projects"; New = "3.1 Safe projects" },
    @{ Style = "Heading 3"; Old = "Scope"; New = "3.1.1 Scope" },
    @{ Style = "Heading 3"; Old = "Examples"; New = "3.1.2 Examples" },
    @{ Style = "Heading 2"; Old = "Process"; New = "3.2 Process" },
    @{ Style = "Heading 3"; Old = "Amendment process"; New = "3.2.1 Amendment process" },
    @{ Style = "Heading 3"; Old = "Appeals process"; New = "3.2.2 Appeals process" }
)

$applied = 0

foreach ($p in $d.Paragraphs) {
    $styleName = $p.Style.NameLocal
    $text = $p.Range.Text

    foreach ($r in $renames) {
        if ($styleName -eq $r.Style -and $text -eq ($r.Old + "`r")) {
            $p.Range.Text = $r.New
            $applied = $applied + 1
            break
        }
    }
}

Write-Output ("Headings renumbered: " + $applied)
